# Applies two corrections to the document:
#  1. "... em boas condições mas não significa ..." -> adds a comma after
#     "condições" and clears the (now stale) grammar-check wavy-underline
#     markers (w:proofErr gramStart/gramEnd) that used to bracket the word.
#  2. Fixes a split word: the paragraph "<h3 id="text-imc">" had the "t"
#     of "text" stranded in the previous run ('<h3 id="t' + 'ext-imc');
#     move it into the run that actually spells the word, without
#     disturbing the spell-check markers around it.

$d = $word.ActiveDocument

# --- Edit 1: "condições" -> "condições," -------------------------------
$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*esquecer dela*") {
        $target1 = $p
        break
    }
}
if ($target1 -ne $null) {
    $r1 = $target1.Range
    # Span the match across both proofErr boundaries (the run that used to
    # be wrapped in gramStart/gramEnd) so the stale grammar markers are
    # dropped along with the text fix, same as Word does when you edit a
    # flagged phrase.
    $r1.Find.Execute("em boas condições mas", $false, $false, $false, `
        $false, $false, $true, 1, $false, "em boas condições, mas", 2) | Out-Null
}

# --- Edit 2: '<h3 id="t' / 'ext-imc' -> '<h3 id="' / 'text-imc' --------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like '*<h3 id="t*' -and $p.Range.Text -notlike '*<h3 id="text-imc*') {
        $target2 = $p
        break
    }
}
if ($target2 -eq $null) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like '*<h3 id="text-imc*') {
            $nextPara = $p.Next()
            if ($nextPara -ne $null -and $nextPara.Range.Text -like "*Mas o que*IMC*") {
                $target2 = $p
                break
            }
        }
    }
}
if ($target2 -ne $null) {
    $r2 = $target2.Range
    $start = $r2.Start
    $end = $r2.End
    $tPos = -1
    for ($i = $start; $i -lt $end; $i++) {
        $c = $d.Range($i, $i + 1)
        if ($c.Text -eq '"') {
            $tPos = $i + 1
            break
        }
    }
    if ($tPos -ge 0) {
        $stray = $d.Range($tPos, $tPos + 1)
        if ($stray.Text -eq "t") {
            # Grow the run that currently reads "ext-imc" so it reads
            # "text-imc" first (this keeps it attached to that run, inside
            # the existing spellStart/spellEnd pair), then delete the
            # orphaned "t" left behind in the preceding run.
            $wordRun = $d.Range($tPos + 1, $tPos + 8)
            $wordRun.Text = "text-imc"
            $stray2 = $d.Range($tPos, $tPos + 1)
            $stray2.Text = ""
        }
    }
}
